$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (July) - H8 changes from 153 to 152
$ws.Range("H8").Value = 152

# Row 9 (August) - label and several values change
$ws.Range("A9").Value = "August (through 08-21)"
$ws.Range("C9").Value = 49
$ws.Range("D9").Value = 55
$ws.Range("E9").Value = 32
$ws.Range("G9").Value = 128
$ws.Range("H9").Value = 106

# Row 10 (Total) - recalculated totals
$ws.Range("C10").Value = 351
$ws.Range("D10").Value = 520
$ws.Range("E10").Value = 457
$ws.Range("G10").Value = 749
$ws.Range("H10").Value = 1020
